$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of an existing "index" cell (column A, style s="1") down to
# the two brand-new rows (16 and 17) before writing their values, so the
# A16/A17 cells keep the same bold/bordered/centered style as the rest of
# column A instead of picking up Excel's unstyled default.
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(17, 1))

# Rows 8-17 (A,B,C,D,E) after inserting line7/line8 before the extr* rows
# and shifting the former extr1..extr8 rows (8..15) down to rows 10..17.
$data = @(
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4", 7,  8,  $false),
    @(14, 12, "extr5", 9,  11, $false),
    @(15, 13, "extr6", 7,  11, $true),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
